$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at K (before existing K:L "Memo Line.." etc.)
[void]$ws.Range("K1:L1").EntireColumn.Insert()
$ws.Range("K1:L1").ColumnWidth = 13

$ws.Range("K1").Value = "Transducer Model"
$ws.Range("K2").Value = "XRS-5"
$ws.Range("L1").Value = "Transducer Serial Number"
$ws.Range("L2").Value = "2007/234500"

# Match the author's final view state: scrolled right with L3 selected
$excel.ActiveWindow.ScrollColumn = 7
$excel.ActiveWindow.ScrollRow = 1
[void]$ws.Range("L3").Select()
